# "Generate Report for Handoff"
#
# The localization-status report was regenerated: the four markdown files
# that are "Ready for handoff" (0a619349…, 4870f5ab…, 9cdfc1a1…, fbc5171f…)
# got a new handoff pass. As a result:
#   - Their Priority moved from "low" to "ht" (on both the zh-cn and
#     de-de per-language sheets).
#   - The "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" for
#     those rows was refreshed to the new generation timestamp
#     (different timestamps per language: zh-cn vs de-de/Overview).

$wb = $excel.ActiveWorkbook

# --- Overview sheet: Latest HO Xliff Generate Date for the 4 handoff rows ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G4:G7").Value = "2016-08-30 12:39:02"

# --- zh-cn sheet: Priority + Latest Handoff Datetime for the 4 handoff rows ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E4:E7").Value = "ht"
$zhcn.Range("H4:H7").Value = "2016-08-30 12:38:56"

# --- de-de sheet: Priority + Latest Handoff Datetime for the 4 handoff rows ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E4:E7").Value = "ht"
$dede.Range("H4:H7").Value = "2016-08-30 12:39:02"
